# Updated cryptos list on Fri Apr 28 08:56:51 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.301.29"
$ws.Range("E2").Value = "  +0.97%  "

# Row 3
$ws.Range("D3").Value = "1.912.99"
$ws.Range("E3").Value = "  +1.39%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.47%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4727"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.92%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4073"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.27%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08032"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.004"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.34%  "

# Row 13
$ws.Range("D13").Value = "1.927.07"
$ws.Range("E13").Value = "  +3.30%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.887"
$ws.Range("D14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.133"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.05%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.00%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06633"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.17%  "

# Row 19
$ws.Range("E19").Value = "  +0.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.50%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "

# Row 22
$ws.Range("D22").Value = "29.326.67"
$ws.Range("E22").Value = "  +0.95%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.510"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.202"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.32%  "

# Row 26
$ws.Range("D26").Value = "2.129.46"
$ws.Range("E26").Value = "  +1.64%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.42%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.77%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.056"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.86%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.111"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.56%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.072"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.12%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09541"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.36%  "

# Row 34
$ws.Range("E34").Value = "  +1.52%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.548"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.48%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.401"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.41%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02254"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.49%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06084"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.73%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.272"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.174"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.47%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5884"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.544"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.96%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1835"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.96%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07876"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.11%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.258"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.93%  "

# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5537"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.51%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.928"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.56%  "

